$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Null-out bank account numbers (column D) for the bank journal rows
$ws.Range("D2").Value = "\N"
$ws.Range("D3").Value = "\N"
$ws.Range("D4").Value = "\N"
$ws.Range("D5").Value = "\N"
$ws.Range("D6").Value = "\N"

# Rename sale journal identifiers
$ws.Range("A8").Value = "z0bug.jou_fatt|inv"
$ws.Range("G8").Value = "FATT"

# Rename purchase journal identifiers
$ws.Range("A9").Value = "z0bug.jou_acq|bill"
$ws.Range("G9").Value = "ACQ"
